# Afshan branch v01 - update lookup codes on the "data" sheet.
# Column B holds short codes; several of them are replaced with new codes.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "PN"    # was "XE"
$ws.Range("B3").Value  = "CF"    # was "H6"
$ws.Range("B5").Value  = "UR"    # was "JT"
$ws.Range("B6").Value  = "0O"    # was "LR"
$ws.Range("B7").Value  = "NGJ"   # was "Q4S"
$ws.Range("B8").Value  = "A8"    # was "UU"
$ws.Range("B13").Value = "0Z"    # was "E5"
